$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 19

$ws.Range("F3").Value = 7
$ws.Range("H3").Value = 7

$ws.Range("E7").Value = 19

$ws.Range("E10").Value = 17

$ws.Range("F13").Value = 3
$ws.Range("H13").Value = 3

$ws.Range("E15").Value = 67
$ws.Range("F15").Value = 35
$ws.Range("H15").Value = 35

$ws.Range("F16").Value = 65
$ws.Range("H16").Value = 65

$ws.Range("E18").Value = 65
$ws.Range("F18").Value = 26
$ws.Range("H18").Value = 26
